$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.651.91"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").Value = "2.498.71"
$ws.Range("E3").Value = "  +9.33%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'480.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.73%  "
$ws.Range("D6").Value = "'138.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.88%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.63%  "
$ws.Range("D9").Value = "2.493.64"
$ws.Range("E9").Value = "  +8.96%  "
$ws.Range("D10").Value = "'0.0988"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.17%  "
$ws.Range("D11").Value = "'5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "'0.326"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.97%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "2.934.19"
$ws.Range("E14").Value = "  +9.36%  "
$ws.Range("D15").Value = "55.753.93"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "'0.0000138"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +15.41%  "
$ws.Range("D17").Value = "'20.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.74%  "
$ws.Range("D18").Value = "2.507.57"
$ws.Range("E18").Value = "  +7.62%  "
$ws.Range("E19").Value = "  +7.10%  "
$ws.Range("D20").Value = "'320.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.44%  "
$ws.Range("D21").Value = "'9.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.15%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'5.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.20%  "
$ws.Range("D24").Value = "'57.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("E27").Value = "  +9.34%  "
$ws.Range("D28").Value = "2.614.38"
$ws.Range("E28").Value = "  +9.18%  "
$ws.Range("D29").Value = "'7.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.51%  "
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("E30").Value = "  +9.61%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "'149.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "'18.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.25%  "
$ws.Range("E34").Value = "  +10.45%  "
$ws.Range("E35").Value = "  +10.23%  "
$ws.Range("D36").Value = "'3.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("E37").Value = "  +10.23%  "
$ws.Range("D38").Value = "'0.841"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'34.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("D40").Value = "'0.613"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +19.55%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0547"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.74%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.43%  "
$ws.Range("D44").Value = "'3.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.75%  "
$ws.Range("D45").Value = "'10.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").Value = "1.972.43"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("E47").Value = "  +11.06%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0223"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.38%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'250.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +31.39%  "
$ws.Range("D50").Value = "'17.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.06%  "
$ws.Range("E51").Value = "  +9.45%  "
